# Update Name of Algo
# Apply the specific numeric corrections to the result data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 6.127499999999989
$ws.Range("E4").Value = 12.7075
$ws.Range("A9").Value = -20.08099999999997
$ws.Range("E10").Value = 12.0886
$ws.Range("A18").Value = -22.90180000000001
$ws.Range("A20").Value = -22.14390000000003
$ws.Range("D21").Value = -7.624200000000002
